# -----------------------------------------------------------------------
# [ADDITIONAL SCRAPING] add code to scrape more data about a player's
# batting performance in a match, also update the existing Excel sheets:
#   * new "Player Info" sheet (first tab)
#   * "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (store bare match code,
#     not the full URL); drop the stray empty INNING_NUMBER cells
#   * "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (store bare match code)
#   * new "ODI Batting Extra" sheet (last tab)
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Set-TextCell($ws, $row, $col, $text) {
    # Force literal text even when the value looks numeric (match codes,
    # "0"/"1" counters, percentages, ...) by using the apostrophe escape,
    # same as typing '<text> into Excel.
    $cell = $ws.Cells.Item($row, $col)
    if ($null -eq $text) {
        $cell.Value = "'"
    } else {
        $cell.Value = "'" + $text
    }
}

# =========================================================================
# 1. New "Player Info" sheet, inserted as the first tab
# =========================================================================
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

Set-HeaderCell $playerInfo 1 1 "ID"
Set-HeaderCell $playerInfo 1 2 "NAME"
Set-HeaderCell $playerInfo 1 3 "BATTING_HAND"
Set-HeaderCell $playerInfo 1 4 "BOWL_STYLE"

Set-TextCell $playerInfo 2 1 "4474"
Set-TextCell $playerInfo 2 2 "Jeffrey Dexter Francis Vandersay"
Set-TextCell $playerInfo 2 3 "Right Handed"
Set-TextCell $playerInfo 2 4 "Right Arm Leg Break"

# =========================================================================
# 2. "ODI Batting": MATCH_CARD_LINK column -> MATCH_CODE column
# =========================================================================
$battingSheet = $wb.Worksheets.Item("ODI Batting")
Set-TextCell $battingSheet 1 4 "MATCH_CODE"

$battingCodes = @{
    2  = "3866"
    3  = "3868"
    4  = "3870"
    5  = "3962"
    6  = "3965"
    7  = "3983"
    8  = "3990"
    9  = "4079"
    10 = "4081"
    11 = "4082"
    12 = "4087"
    13 = "4344"
    14 = "4521"
    15 = "4523"
    16 = "4527"
    17 = "4597"
    18 = "4600"
    19 = "4601"
    20 = "4603"
    21 = "4691"
}
foreach ($row in $battingCodes.Keys) {
    Set-TextCell $battingSheet $row 4 $battingCodes[$row]
}

# Rows whose INNING_NUMBER (column B) cell was a stray empty string - drop
# them entirely (genuinely blank, not an empty-text cell).
$emptyInningRows = @(3, 4, 5, 6, 13, 14, 18, 21)
foreach ($row in $emptyInningRows) {
    $battingSheet.Cells.Item($row, 2).Value = ""
}

# =========================================================================
# 3. "ODI Bowling": MATCH_CARD_LINK column -> MATCH_CODE column
# =========================================================================
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
Set-TextCell $bowlingSheet 1 2 "MATCH_CODE"

$bowlingCodes = @{
    2  = "3866"
    3  = "3868"
    4  = "3870"
    5  = "3965"
    6  = "3983"
    7  = "3990"
    8  = "4079"
    9  = "4081"
    10 = "4082"
    11 = "4087"
    12 = "4344"
    13 = "4521"
    14 = "4523"
    15 = "4527"
    16 = "4597"
    17 = "4600"
    18 = "4601"
    19 = "4603"
    20 = "4691"
}
foreach ($row in $bowlingCodes.Keys) {
    Set-TextCell $bowlingSheet $row 2 $bowlingCodes[$row]
}

# =========================================================================
# 4. New "ODI Batting Extra" sheet, inserted as the last tab
# =========================================================================
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extraSheet = $wb.Worksheets.Add($null, $bowlingSheet)
$extraSheet.Name = "ODI Batting Extra"

Set-HeaderCell $extraSheet 1 1 "MATCH_CODE"
Set-HeaderCell $extraSheet 1 2 "BATTING_POSITION"
Set-HeaderCell $extraSheet 1 3 "NUM_4"
Set-HeaderCell $extraSheet 1 4 "NUM_6"
Set-HeaderCell $extraSheet 1 5 "PERCENT_RUNS_OF_TOTAL"
Set-HeaderCell $extraSheet 1 6 "MAN_OF_MATCH"

$extraRows = @(
    @{A="3866"; B=11;     C="0"; D="0"; E="5.98%";  F="NO"},
    @{A="3868"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="3870"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="3962"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="3965"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="3983"; B=9;      C="1"; D="0"; E="4.42%";  F="NO"},
    @{A="3990"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="4079"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="4081"; B=9;      C="0"; D="0"; E="11.76%"; F="NO"},
    @{A="4082"; B=7;      C="0"; D="0"; E=$null;    F="NO"},
    @{A="4087"; B=9;      C="0"; D="0"; E="1.94%";  F="NO"},
    @{A="4344"; B=9;      C=$null; D=$null; E=$null;    F="NO"},
    @{A="4521"; B=8;      C=$null; D=$null; E=$null;    F="NO"},
    @{A="4523"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="4527"; B=10;     C="0"; D="0"; E="1.18%";  F="NO"},
    @{A="4597"; B=9;      C="0"; D="0"; E="3.18%";  F="NO"},
    @{A="4600"; B=$null;  C=$null; D=$null; E=$null;    F="NO"},
    @{A="4601"; B=10;     C="0"; D="0"; E=$null;    F="NO"},
    @{A="4603"; B=9;      C="0"; D="0"; E="2.50%";  F="NO"},
    @{A="4691"; B=$null;  C=$null; D=$null; E=$null;    F="NO"}
)

$row = 2
foreach ($data in $extraRows) {
    Set-TextCell $extraSheet $row 1 $data.A

    if ($null -eq $data.B) {
        $extraSheet.Cells.Item($row, 2).Value = "'"
    } else {
        $extraSheet.Cells.Item($row, 2).Value = $data.B
    }

    Set-TextCell $extraSheet $row 3 $data.C
    Set-TextCell $extraSheet $row 4 $data.D
    Set-TextCell $extraSheet $row 5 $data.E
    Set-TextCell $extraSheet $row 6 $data.F

    $row++
}
